$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 becomes the Puerto Rico - Argentina match (final score 0:6, picked & correct)
$ws.Range("A2").Value = "Puerto Rico - Argentina ✓: 0:6"
$ws.Range("B2").Value = "Argentina"
$ws.Range("C2").Value = 68
$ws.Range("D2").Value = 50
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = 1.7
$ws.Range("G2").Value = "✓"

# Row 3 becomes the Palmeiras - Red Bull Bragantino match (final score 5:1, picked & correct)
$ws.Range("A3").Value = "Sociedade Esportiva Palmeiras ✓ - Red Bull Bragantino: 5:1"
$ws.Range("B3").Value = "Sociedade Esportiva Palmeiras"
$ws.Range("C3").Value = 68
$ws.Range("D3").Value = 92
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = 1.4
$ws.Range("G3").Value = "✓"
